$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 697, shifting existing rows (697-738) down to (698-739)
$ws.Rows("697:697").Insert()

# Populate the newly inserted row with the new data point.
# Column A holds the date as plain text (matches the rest of the sheet), so
# briefly force a text number format before writing the value to stop Excel
# from auto-converting the "yyyy/mm/dd"-looking string into a real date
# serial, then restore the cell's style to match its neighbours (plain,
# unstyled data cells) so only the value itself changed.
$ws.Cells.Item(697, 1).NumberFormat = "@"
$ws.Cells.Item(697, 1).Value = "2026/01/23"
$ws.Cells.Item(697, 1).Style = $ws.Cells.Item(696, 1).Style
$ws.Cells.Item(697, 2).Value = "金"
$ws.Cells.Item(697, 3).Value = 7
$ws.Cells.Item(697, 4).Value = 201
